$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Addr, $Text) {
    $rng = $Sheet.Range($Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

Set-TextValue $ws "D2" "306.80"
Set-TextValue $ws "E2" "1.60%"
Set-TextValue $ws "D3" "35.91"
Set-TextValue $ws "E3" "1.20%"
Set-TextValue $ws "D4" "5.055"
Set-TextValue $ws "E4" "0.07%"
Set-TextValue $ws "D5" "0.08087"
Set-TextValue $ws "E5" "1.05%"
Set-TextValue $ws "D6" "1.945"
Set-TextValue $ws "E6" "0.60%"
Set-TextValue $ws "D7" "4.149"
Set-TextValue $ws "E7" "2.25%"
Set-TextValue $ws "D8" "7.775"
Set-TextValue $ws "E8" "-0.45%"
Set-TextValue $ws "D9" "0.9302"
Set-TextValue $ws "E9" "0.52%"
Set-TextValue $ws "D10" "0.1356"
Set-TextValue $ws "E10" "1.03%"
Set-TextValue $ws "D11" "0.1905"
Set-TextValue $ws "E11" "0.50%"
Set-TextValue $ws "D12" "0.09247"
Set-TextValue $ws "E12" "0.13%"
Set-TextValue $ws "D13" "0.03522"
Set-TextValue $ws "E13" "4.05%"
Set-TextValue $ws "D14" "0.09864"
Set-TextValue $ws "E14" "-0.03%"
Set-TextValue $ws "D15" "0.001434"
Set-TextValue $ws "E15" "2.83%"
Set-TextValue $ws "D16" "0.005812"
Set-TextValue $ws "E16" "0.86%"
Set-TextValue $ws "D17" "3.571"
Set-TextValue $ws "E17" "1.98%"
Set-TextValue $ws "D19" "0.3445"
Set-TextValue $ws "E19" "1.22%"
Set-TextValue $ws "D20" "0.1345"
Set-TextValue $ws "E20" "3.24%"
Set-TextValue $ws "D21" "4.886"
Set-TextValue $ws "E21" "-3.59%"
Set-TextValue $ws "D22" "0.2596"
Set-TextValue $ws "E22" "8.08%"
Set-TextValue $ws "D23" "0.04388"
Set-TextValue $ws "E23" "-2.39%"
Set-TextValue $ws "D24" "0.001221"
Set-TextValue $ws "E24" "0.49%"
Set-TextValue $ws "D25" "0.004776"
Set-TextValue $ws "E25" "-0.47%"
Set-TextValue $ws "E26" "31.87%"
Set-TextValue $ws "D27" "0.0003125"
Set-TextValue $ws "E27" "3.98%"
Set-TextValue $ws "D39" "0.01982"
Set-TextValue $ws "E39" "3.89%"
Set-TextValue $ws "D40" "0.05077"
Set-TextValue $ws "E40" "7.21%"
Set-TextValue $ws "D41" "0.01121"
Set-TextValue $ws "E41" "16.05%"
Set-TextValue $ws "D42" "0.007612"
Set-TextValue $ws "E42" "3.28%"
Set-TextValue $ws "D43" "0.1376"
Set-TextValue $ws "E43" "3.17%"
Set-TextValue $ws "D44" "0.002097"
Set-TextValue $ws "E44" "-0.76%"
Set-TextValue $ws "D45" "0.01082"
Set-TextValue $ws "E45" "2.70%"
Set-TextValue $ws "D46" "0.00006389"
Set-TextValue $ws "E46" "0.65%"
Set-TextValue $ws "D47" "0.00000000749"
Set-TextValue $ws "E47" "-0.25%"
Set-TextValue $ws "D48" "65.22"
Set-TextValue $ws "E48" "1.15%"
Set-TextValue $ws "D49" "0.001189"
Set-TextValue $ws "E49" "-28.44%"
Set-TextValue $ws "D50" "0.00002098"
Set-TextValue $ws "E50" "-0.25%"
Set-TextValue $ws "D51" "0.0001998"
Set-TextValue $ws "E51" "-0.25%"
